$wb = $excel.ActiveWorkbook

# --- Update conversion text on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.34 = 13014.02 pesos`n✅ 13014.02 pesos = 3.33 = 970.23 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update rate values on "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("O10").Value = 3891.05
$ws2.Range("O12").Value = 291.5
